$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: same style as E1 (bold header style)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Data cells F2:F12 - timestamps, stored as text (no special style, matches column E data cells)
$timestamps = @(
    "2021-10-05 13:39:19.183268",
    "2021-10-05 13:39:19.183280",
    "2021-10-05 13:39:19.183284",
    "2021-10-05 13:39:19.183287",
    "2021-10-05 13:39:19.183291",
    "2021-10-05 13:39:19.183294",
    "2021-10-05 13:39:19.183297",
    "2021-10-05 13:39:19.183300",
    "2021-10-05 13:39:19.183303",
    "2021-10-05 13:39:19.183306",
    "2021-10-05 13:39:19.183309"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
